$wb = $excel.ActiveWorkbook

# Sheet ALC row 13
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1192.5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1192.5
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1192.5
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1530.5

# Sheet ALC row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 72593.8
$ws.Range("I47").Value = 58000
$ws.Range("J47").Value = 82323
$ws.Range("K47").Value = 58000
$ws.Range("L47").Value = 82323
$ws.Range("M47").Value = -57028
$ws.Range("N47").Value = -84267

# Sheet ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2790

# Sheet ALC row 54
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 73333.336
$ws.Range("I54").Value = 68000
$ws.Range("J54").Value = 84000
$ws.Range("K54").Value = 68000
$ws.Range("L54").Value = 84000
$ws.Range("M54").Value = -67514
$ws.Range("N54").Value = -84972

# Sheet ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2332
$ws.Range("I113").Value = 2498
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2498
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 756
$ws.Range("N113").Value = -8508

# Sheet ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2150.2942
$ws.Range("I2").Value = 1672.75
$ws.Range("J2").Value = 3296.4
$ws.Range("K2").Value = 1672.75
$ws.Range("L2").Value = 3296.4
$ws.Range("M2").Value = -1559.75
$ws.Range("N2").Value = -3522.4

# Sheet ARM row 14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 456
$ws.Range("I14").Value = 456
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 456
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -281
$ws.Range("N14").ClearContents()

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5358.71
$ws.Range("I32").Value = 2971.7114
$ws.Range("J32").Value = 12660.117
$ws.Range("K32").Value = 2971.7114
$ws.Range("L32").Value = 12660.117
$ws.Range("M32").Value = -2684.7114
$ws.Range("N32").Value = -13234.117

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12608527
$ws.Range("I45").Value = 2252.4
$ws.Range("K45").Value = 2252.4
$ws.Range("M45").Value = -1875.4

# Sheet ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2150.2942
$ws.Range("I116").Value = 1672.75
$ws.Range("J116").Value = 3296.4
$ws.Range("K116").Value = 1672.75
$ws.Range("L116").Value = 3296.4
$ws.Range("M116").Value = 621.25
$ws.Range("N116").Value = -7884.4

# Sheet BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2150.2942
$ws.Range("I3").Value = 1672.75
$ws.Range("J3").Value = 3296.4
$ws.Range("K3").Value = 1672.75
$ws.Range("L3").Value = 3296.4
$ws.Range("M3").Value = -1558.75
$ws.Range("N3").Value = -3524.4

# Sheet BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 111544.22
$ws.Range("I22").Value = 143292.58
$ws.Range("J22").Value = 425
$ws.Range("K22").Value = 143292.58
$ws.Range("L22").Value = 425
$ws.Range("M22").Value = -143119.58
$ws.Range("N22").Value = -771

# Sheet BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3794.56
$ws.Range("I107").Value = 4623.5713
$ws.Range("J107").Value = 3472.1667
$ws.Range("K107").Value = 4623.5713
$ws.Range("L107").Value = 3472.1667
$ws.Range("M107").Value = -2703.5713
$ws.Range("N107").Value = -7312.1667

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5079.826
$ws.Range("I134").Value = 3545.8
$ws.Range("K134").Value = 10637.4
$ws.Range("M134").Value = -8102.400000000001

# Sheet CRP row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 4993
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 4993
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4993
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -5333

# Sheet CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3096.9644
$ws.Range("I105").Value = 1838.375
$ws.Range("J105").Value = 4775.0835
$ws.Range("K105").Value = 1838.375
$ws.Range("L105").Value = 4775.0835
$ws.Range("M105").Value = -91.375
$ws.Range("N105").Value = -8269.083500000001

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2816.5186
$ws.Range("I122").Value = 2146.8462
$ws.Range("J122").Value = 3438.3572
$ws.Range("K122").Value = 6440.5386
$ws.Range("L122").Value = 10315.0716
$ws.Range("M122").Value = -3990.5386
$ws.Range("N122").Value = -15215.0716

# Sheet CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2581.6667
$ws.Range("I68").Value = 960
$ws.Range("K68").Value = 2880
$ws.Range("M68").Value = -2069

# Sheet CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2581.6667
$ws.Range("I71").Value = 960
$ws.Range("K71").Value = 8640
$ws.Range("M71").Value = -4584

# Sheet CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1040.0714
$ws.Range("J107").Value = 1028
$ws.Range("L107").Value = 3084
$ws.Range("N107").Value = -6924

# Sheet CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Sheet CUL row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 454997.5
$ws.Range("I128").Value = 454997.5
$ws.Range("K128").Value = 1364992.5
$ws.Range("M128").Value = -1360012.5

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1304.2106
$ws.Range("I131").Value = 887.63635
$ws.Range("J131").Value = 1877
$ws.Range("K131").Value = 2662.90905
$ws.Range("L131").Value = 5631
$ws.Range("M131").Value = 2377.09095
$ws.Range("N131").Value = -15711

# Sheet CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1251.9286
$ws.Range("I139").Value = 1251.9286
$ws.Range("K139").Value = 3755.7858
$ws.Range("M139").Value = 1384.2142

# Sheet GSM row 9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 200
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -30
$ws.Range("N9").ClearContents()

# Sheet GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 27299.4
$ws.Range("J57").Value = 30499.25
$ws.Range("L57").Value = 30499.25
$ws.Range("N57").Value = -32139.25

# Sheet GSM row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 22750
$ws.Range("J63").Value = 22750
$ws.Range("L63").Value = 22750
$ws.Range("N63").Value = -24122

# Sheet GSM row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 22750
$ws.Range("J66").Value = 22750
$ws.Range("L66").Value = 68250
$ws.Range("N66").Value = -75114

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7018
$ws.Range("I122").Value = 9242.6
$ws.Range("J122").Value = 4237.25
$ws.Range("K122").Value = 27727.8
$ws.Range("L122").Value = 12711.75
$ws.Range("M122").Value = -25277.8
$ws.Range("N122").Value = -17611.75

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11199.75
$ws.Range("I132").Value = 9194.6
$ws.Range("J132").Value = 14541.667
$ws.Range("K132").Value = 27583.8
$ws.Range("L132").Value = 43625.001
$ws.Range("M132").Value = -25053.8
$ws.Range("N132").Value = -48685.001

# Sheet LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 6324.3184
$ws.Range("I55").Value = 1922.3334
$ws.Range("J55").Value = 11606.7
$ws.Range("K55").Value = 1922.3334
$ws.Range("L55").Value = 11606.7
$ws.Range("M55").Value = -1749.3334
$ws.Range("N55").Value = -11952.7

# Sheet LTW row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 31061.834
$ws.Range("J106").Value = 31061.834
$ws.Range("L106").Value = 31061.834
$ws.Range("N106").Value = -33585.834

# Sheet LTW row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 62629652
$ws.Range("I122").Value = 71575530
$ws.Range("J122").Value = 8505
$ws.Range("K122").Value = 214726590
$ws.Range("L122").Value = 25515
$ws.Range("M122").Value = -214724140
$ws.Range("N122").Value = -30415

# Sheet LTW row 124
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 89999.5
$ws.Range("J124").Value = 89999.5
$ws.Range("L124").Value = 89999.5
$ws.Range("N124").Value = -99819.5

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4426.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4426.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13278.75
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -18338.75

# Sheet WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4000
$ws.Range("I14").Value = 4000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -3832
$ws.Range("N14").ClearContents()

# Sheet WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13737.625
$ws.Range("J62").Value = 13500
$ws.Range("L62").Value = 13500
$ws.Range("N62").Value = -14748

# Sheet WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 13737.625
$ws.Range("J65").Value = 13500
$ws.Range("L65").Value = 67500
$ws.Range("N65").Value = -73740
